$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, [string]$value) {
    # Build a formula that evaluates to the literal text, then convert to a plain value
    # via copy / paste-special so the stored cell keeps type "string" without altering styles.
    $escaped = $value.Replace('"', '""')
    $cellRange.Formula = '="' + $escaped + '"'
    $cellRange.Copy()
    $cellRange.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "66.430.62"
Set-TextValue $ws.Range("E2") "  +0.02%  "
Set-TextValue $ws.Range("D3") "3.216.00"
Set-TextValue $ws.Range("E3") "  +0.91%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "608.78"
Set-TextValue $ws.Range("E5") "  +2.09%  "
Set-TextValue $ws.Range("D6") "157.23"
Set-TextValue $ws.Range("E6") "  +1.94%  "
Set-TextValue $ws.Range("D8") "3.216.49"
Set-TextValue $ws.Range("E8") "  +1.05%  "
Set-TextValue $ws.Range("E9") "  +0.73%  "
Set-TextValue $ws.Range("E10") "  +0.61%  "
Set-TextValue $ws.Range("D11") "5.71"
Set-TextValue $ws.Range("E11") "  -4.12%  "
Set-TextValue $ws.Range("E12") "  -2.58%  "
Set-TextValue $ws.Range("E13") "  +0.79%  "
Set-TextValue $ws.Range("D14") "38.68"
Set-TextValue $ws.Range("E14") "  -1.29%  "
Set-TextValue $ws.Range("D15") "3.744.97"
Set-TextValue $ws.Range("E15") "  +1.08%  "
Set-TextValue $ws.Range("D16") "66.505.75"
Set-TextValue $ws.Range("E16") "  +0.24%  "
Set-TextValue $ws.Range("D17") "7.37"
Set-TextValue $ws.Range("E17") "  -1.36%  "
Set-TextValue $ws.Range("D18") "3.222.22"
Set-TextValue $ws.Range("E18") "  +1.17%  "
Set-TextValue $ws.Range("E19") "  +1.31%  "
Set-TextValue $ws.Range("D20") "507.47"
Set-TextValue $ws.Range("E20") "  -1.44%  "
Set-TextValue $ws.Range("E21") "  -1.28%  "
Set-TextValue $ws.Range("D22") "0.732"
Set-TextValue $ws.Range("E22") "  -0.71%  "
Set-TextValue $ws.Range("D23") "8.00"
Set-TextValue $ws.Range("E23") "  -0.81%  "
Set-TextValue $ws.Range("D24") "14.63"
Set-TextValue $ws.Range("E24") "  -2.07%  "
Set-TextValue $ws.Range("E25") "  -0.85%  "
Set-TextValue $ws.Range("E26") "  +0.05%  "
Set-TextValue $ws.Range("D27") "3.01"
Set-TextValue $ws.Range("E27") "  +0.15%  "
Set-TextValue $ws.Range("D28") "9.09"
Set-TextValue $ws.Range("E28") "  -2.08%  "
Set-TextValue $ws.Range("E29") "  +1.36%  "
Set-TextValue $ws.Range("D30") "0.126"
Set-TextValue $ws.Range("E30") "  +39.89%  "
Set-TextValue $ws.Range("E31") "  +0.12%  "
Set-TextValue $ws.Range("E32") "  -1.80%  "
Set-TextValue $ws.Range("E33") "  -0.48%  "
Set-TextValue $ws.Range("E34") "  +0.23%  "
Set-TextValue $ws.Range("E35") "  -4.20%  "
Set-TextValue $ws.Range("D36") "6.51"
Set-TextValue $ws.Range("E36") "  -0.20%  "
Set-TextValue $ws.Range("D37") "504.76"
Set-TextValue $ws.Range("E37") "  -0.59%  "
Set-TextValue $ws.Range("D38") "55.33"
Set-TextValue $ws.Range("E38") "  +0.82%  "
Set-TextValue $ws.Range("D39") "0.0₃0774"
Set-TextValue $ws.Range("E39") "  +14.55%  "
Set-TextValue $ws.Range("D40") "3.08"
Set-TextValue $ws.Range("E40") "  +6.90%  "
Set-TextValue $ws.Range("E41") "  -0.74%  "
Set-TextValue $ws.Range("E42") "  +3.13%  "
Set-TextValue $ws.Range("E43") "  -1.85%  "
Set-TextValue $ws.Range("D44") "0.298"
Set-TextValue $ws.Range("E44") "  -1.84%  "
Set-TextValue $ws.Range("D45") "2.46"
Set-TextValue $ws.Range("E45") "  +0.61%  "
Set-TextValue $ws.Range("D46") "2.915.17"
Set-TextValue $ws.Range("E46") "  +0.19%  "
Set-TextValue $ws.Range("D47") "28.24"
Set-TextValue $ws.Range("E47") "  -1.51%  "
Set-TextValue $ws.Range("E48") "  +3.36%  "
Set-TextValue $ws.Range("E50") "  -0.77%  "
Set-TextValue $ws.Range("D51") "122.03"
Set-TextValue $ws.Range("E51") "  -0.89%  "
